$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "BBB"
$ws.Range("B5").Value = "BBB"
$ws.Range("B6").Value = "BB"
$ws.Range("B9").Value = "A"
$ws.Range("B10").Value = "BBB"
$ws.Range("B12").Value = "BB"
$ws.Range("B15").Value = "B"
$ws.Range("B23").Value = "BBB"
$ws.Range("B24").Value = "BBB"
$ws.Range("B28").Value = "BBB"
$ws.Range("B29").Value = "B"
$ws.Range("B38").Value = "BBB"
$ws.Range("B47").Value = "BBB"
$ws.Range("B48").Value = "BBB"
$ws.Range("B50").Value = "BBB"
$ws.Range("B53").Value = "BB"
$ws.Range("B64").Value = "B"
$ws.Range("B65").Value = "BBB"
$ws.Range("B66").Value = "BBB"
$ws.Range("B67").Value = "A"
$ws.Range("B72").Value = "BBB"
$ws.Range("B80").Value = "BBB"
$ws.Range("B87").Value = "BBB"
$ws.Range("B88").Value = "B"
$ws.Range("B90").Value = "BBB"
$ws.Range("B93").Value = "BB"
$ws.Range("B98").Value = "BBB"
$ws.Range("B99").Value = "B"
$ws.Range("B100").Value = "B"
$ws.Range("B106").Value = "BBB"
$ws.Range("B108").Value = "BB"
$ws.Range("B109").Value = "BB"
$ws.Range("B110").Value = "BBB"
$ws.Range("B111").Value = "BB"
$ws.Range("B121").Value = "BB"
$ws.Range("B123").Value = "BB"
$ws.Range("B126").Value = "BB"
$ws.Range("B128").Value = "BBB"
$ws.Range("B129").Value = "BB"
$ws.Range("B130").Value = "BB"
$ws.Range("B135").Value = "AA"
$ws.Range("B140").Value = "BBB"
$ws.Range("B143").Value = "BBB"
$ws.Range("B144").Value = "BBB"
$ws.Range("B148").Value = "AA"
$ws.Range("B151").Value = "BBB"
$ws.Range("B153").Value = "BB"
$ws.Range("B156").Value = "BB"
$ws.Range("B158").Value = "BBB"
$ws.Range("B159").Value = "BB"
$ws.Range("B163").Value = "BB"
$ws.Range("B165").Value = "B"
$ws.Range("B166").Value = "AA"
$ws.Range("B169").Value = "BB"
$ws.Range("B170").Value = "BBB"
$ws.Range("B171").Value = "BBB"
$ws.Range("B179").Value = "BBB"
$ws.Range("B180").Value = "BBB"
$ws.Range("B182").Value = "BB"
$ws.Range("B192").Value = "BBB"
$ws.Range("B195").Value = "BB"
$ws.Range("B196").Value = "BB"
$ws.Range("B199").Value = "BBB"
$ws.Range("B200").Value = "AA"
$ws.Range("B208").Value = "BBB"
$ws.Range("B220").Value = "BBB"
$ws.Range("B221").Value = "A"
$ws.Range("B222").Value = "AA"
$ws.Range("B226").Value = "BB"
$ws.Range("B229").Value = "B"
$ws.Range("B230").Value = "BBB"
$ws.Range("B238").Value = "BBB"
$ws.Range("B253").Value = "BB"
$ws.Range("B254").Value = "BB"
$ws.Range("B255").Value = "BBB"
$ws.Range("B256").Value = "BB"
$ws.Range("B260").Value = "BB"
$ws.Range("B261").Value = "BB"
$ws.Range("B264").Value = "BBB"
$ws.Range("B265").Value = "BB"
$ws.Range("B266").Value = "B"
$ws.Range("B267").Value = "A"
$ws.Range("B273").Value = "AA"
$ws.Range("B274").Value = "BB"
$ws.Range("B275").Value = "BBB"
$ws.Range("B276").Value = "BBB"
$ws.Range("B282").Value = "BB"
$ws.Range("B283").Value = "BB"
$ws.Range("B288").Value = "A"
$ws.Range("B290").Value = "BBB"
$ws.Range("B296").Value = "BBB"
$ws.Range("B297").Value = "BBB"
$ws.Range("B298").Value = "BB"
$ws.Range("B299").Value = "BBB"
$ws.Range("B302").Value = "BB"
$ws.Range("B305").Value = "BB"
$ws.Range("B306").Value = "BB"
$ws.Range("B311").Value = "BBB"
$ws.Range("B317").Value = "BBB"
$ws.Range("B318").Value = "BBB"
$ws.Range("B323").Value = "BB"
$ws.Range("B324").Value = "B"
$ws.Range("B330").Value = "AA"
$ws.Range("B333").Value = "BBB"
$ws.Range("B337").Value = "BBB"
$ws.Range("B339").Value = "BB"
$ws.Range("B340").Value = "BB"
$ws.Range("B344").Value = "BB"
$ws.Range("B349").Value = "BBB"
$ws.Range("B350").Value = "B"
$ws.Range("B353").Value = "BB"
$ws.Range("B354").Value = "BB"
$ws.Range("B356").Value = "A"
$ws.Range("B358").Value = "A"
$ws.Range("B363").Value = "BBB"
$ws.Range("B366").Value = "BB"
$ws.Range("B367").Value = "BB"
$ws.Range("B368").Value = "BBB"
$ws.Range("B371").Value = "BBB"
$ws.Range("B376").Value = "BBB"
$ws.Range("B384").Value = "BB"
$ws.Range("B387").Value = "BBB"
$ws.Range("B388").Value = "BBB"
$ws.Range("B390").Value = "BB"
$ws.Range("B392").Value = "BBB"
$ws.Range("B395").Value = "BB"
$ws.Range("B397").Value = "CCC"
$ws.Range("B398").Value = "BB"
$ws.Range("B399").Value = "BBB"
$ws.Range("B400").Value = "B"
$ws.Range("B401").Value = "AAA"
$ws.Range("B403").Value = "BBB"
